$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2:D7").Sort($ws1.Range("A2:A7"), 1)

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2:L7").Sort($ws2.Range("A2:A7"), 1)

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2:L7").Sort($ws3.Range("A2:A7"), 1)
